$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at 275, shifting existing rows 275:283 down to 276:284
$ws.Range("A275").EntireRow.Insert()

# Populate the newly inserted row 275 with the new weekly record
$ws.Range("A275").Value = 5
$ws.Range("B275").Value = "Macroferia Regional de Talca"
$ws.Range("C275").Value = "Maule"
$ws.Range("D275").Value = 45041
$ws.Range("E275").Value = 7
$ws.Range("F275").Value = 100112017
$ws.Range("G275").Value = "Apio"
$ws.Range("H275").Value = "Americana (o)"
$ws.Range("I275").Value = "Primera"
$ws.Range("J275").Value = 600
$ws.Range("K275").Value = 7000
$ws.Range("L275").Value = 7000
$ws.Range("M275").Value = 7000
$ws.Range("N275").Value = "`$/docena de matas"
$ws.Range("O275").Value = "Provincia del Elquí"
$ws.Range("P275").Value = 1167
$ws.Range("Q275").Value = 6
$ws.Range("R275").Value = "Hortaliza"
